$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.66"
$ws.Range("E2").Value = "'1.15%"
$ws.Range("D3").Value = "'27.13"
$ws.Range("E3").Value = "'0.94%"
$ws.Range("D4").Value = "'4.701"
$ws.Range("E4").Value = "'0.55%"
$ws.Range("D5").Value = "'0.06188"
$ws.Range("E5").Value = "'3.29%"
$ws.Range("D6").Value = "'6.691"
$ws.Range("E6").Value = "'0.44%"
$ws.Range("D7").Value = "'0.8500"
$ws.Range("E7").Value = "'-0.68%"
$ws.Range("D8").Value = "'0.9154"
$ws.Range("E8").Value = "'-0.69%"
$ws.Range("D9").Value = "'0.1407"
$ws.Range("E9").Value = "'1.29%"
$ws.Range("D10").Value = "'0.04622"
$ws.Range("E10").Value = "'-11.84%"
$ws.Range("D11").Value = "'0.07078"
$ws.Range("E11").Value = "'1.09%"
$ws.Range("E12").Value = "'3.56%"
$ws.Range("D13").Value = "'0.09033"
$ws.Range("E13").Value = "'-1.06%"
$ws.Range("D14").Value = "'0.001536"
$ws.Range("E14").Value = "'-0.40%"
$ws.Range("D15").Value = "'0.0006155"
$ws.Range("E15").Value = "'1.85%"
$ws.Range("D16").Value = "'0.006135"
$ws.Range("E16").Value = "'1.40%"
$ws.Range("D17").Value = "'3.459"
$ws.Range("E17").Value = "'0.12%"
$ws.Range("E18").Value = "'0.85%"
$ws.Range("E20").Value = "'-0.97%"
$ws.Range("E21").Value = "'0.88%"
$ws.Range("D22").Value = "'4.113"
$ws.Range("E22").Value = "'-0.84%"
$ws.Range("D23").Value = "'0.04223"
$ws.Range("E23").Value = "'-0.54%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-0.06%"
$ws.Range("E25").Value = "'-5.86%"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E27").Value = "'-6.51%"
$ws.Range("D40").Value = "'0.03917"
$ws.Range("E40").Value = "'2.06%"
$ws.Range("E41").Value = "'-0.25%"
$ws.Range("D42").Value = "'0.004106"
$ws.Range("E42").Value = "'7.96%"
$ws.Range("E43").Value = "'-9.73%"
$ws.Range("E44").Value = "'-7.68%"
$ws.Range("D45").Value = "'0.00005134"
$ws.Range("E45").Value = "'0.76%"
$ws.Range("E46").Value = "'0.13%"
$ws.Range("D48").Value = "'0.1667"
$ws.Range("E48").Value = "'10.88%"
$ws.Range("E49").Value = "'0.13%"
$ws.Range("E50").Value = "'0.13%"
